# Actualización al 10 de junio
$wb = $excel.ActiveWorkbook

# --- Sheet "Ingreso" (Ingresos/aportes) ---
$wsIngreso = $wb.Worksheets.Item("Ingreso")

$ingresoRows = @(
    @(45081, "Anuel",   100),
    @(45081, "Carlos",  200),
    @(45081, "Randy",   100),
    @(45081, "Julio",   100),
    @(45081, "kukito",   50),
    @(45081, "Punto",   400),
    @(45081, "Alfredo", 100),
    @(45081, "Johan",   300),
    @(45074, "Wilkin",  100),
    @(45074, "Omaury",  100),
    @(45074, "Jeicol",  100),
    @(45074, "Anuel",    85)
)

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$startRow = 432
for ($i = 0; $i -lt $ingresoRows.Count; $i++) {
    $r = $startRow + $i
    $data = $ingresoRows[$i]
    $date = $epoch.AddDays([double]$data[0])
    $wsIngreso.Cells.Item($r, 1).Value = $date
    $wsIngreso.Cells.Item($r, 2).Value = $data[1]
    $wsIngreso.Cells.Item($r, 3).Value = $data[2]
    $wsIngreso.Cells.Item($r, 4).Value = "Aporte"
    if ($r -le 438) {
        $wsIngreso.Cells.Item($r, 3).Style = "Normal"
    }
}

# --- Sheet "Gastos" ---
$wsGastos = $wb.Worksheets.Item("Gastos")

$gastosDate = $epoch.AddDays(45081)
$wsGastos.Cells.Item(44, 1).Value = $gastosDate
$wsGastos.Cells.Item(44, 2).Value = "Agua"
$wsGastos.Cells.Item(44, 3).Value = 150

# --- Update view/selection state to match latest edits ---
$wsGastos.Activate() | Out-Null
$wsGastos.Range("A44").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1

$wsIngreso.Activate() | Out-Null
$wsIngreso.Range("A441:A443").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 423
$excel.ActiveWindow.ScrollColumn = 1
